# Apply updated crypto price/volume values to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel;
# format them as Text first so the exact string (incl. trailing zeros) survives.
$textFormatCells = @("D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D18","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    'D2' = '26.613.89'
    'E2' = '  -0.09%  '
    'D3' = '1.849.78'
    'E3' = '  -0.25%  '
    'E4' = '  +0.11%  '
    'D5' = '263.10'
    'E5' = '  -0.70%  '
    'E6' = '  +0.15%  '
    'D7' = '0.5347'
    'E7' = '  +2.23%  '
    'D8' = '0.3151'
    'E8' = '  -4.04%  '
    'D9' = '0.06931'
    'E9' = '  +1.93%  '
    'D10' = '18.85'
    'E10' = '  -0.07%  '
    'D11' = '0.7706'
    'E11' = '  -0.80%  '
    'D12' = '0.07839'
    'E12' = '  +1.15%  '
    'D13' = '1.852.91'
    'E13' = '  -0.75%  '
    'D14' = '89.73'
    'E14' = '  +1.22%  '
    'D15' = '5.053'
    'E15' = '  +0.39%  '
    'D16' = '14.11'
    'E16' = '  +0.68%  '
    'E17' = '  +0.08%  '
    'D18' = '0.000007963'
    'E18' = '  -0.22%  '
    'E19' = '  +0.11%  '
    'D20' = '26.646.74'
    'E20' = '  -0.04%  '
    'D21' = '2.089.03'
    'E21' = '  +0.24%  '
    'D22' = '4.653'
    'E22' = '  +0.25%  '
    'D23' = '6.024'
    'E23' = '  +0.28%  '
    'D24' = '9.355'
    'E24' = '  -2.03%  '
    'D25' = '2.216'
    'E25' = '  +0.73%  '
    'D26' = '141.52'
    'E26' = '  -1.97%  '
    'E27' = '  +0.92%  '
    'D28' = '17.12'
    'E28' = '  +0.70%  '
    'D29' = '111.76'
    'D30' = '4.306'
    'E30' = '  +2.73%  '
    'D31' = '0.08780'
    'E31' = '  +0.18%  '
    'D32' = '4.115'
    'E32' = '  -1.00%  '
    'D33' = '0.04857'
    'E33' = '  +0.52%  '
    'D34' = '0.7380'
    'E34' = '  +3.16%  '
    'D35' = '1.139'
    'E35' = '  +0.12%  '
    'D36' = '2.905'
    'E36' = '  +1.18%  '
    'D37' = '3.115'
    'E37' = '  -0.05%  '
    'D38' = '2.341'
    'E38' = '  +6.58%  '
    'E39' = '  -2.80%  '
    'D40' = '0.4831'
    'E40' = '  -1.18%  '
    'D41' = '0.9064'
    'E41' = '  +0.53%  '
    'D42' = '108.51'
    'E42' = '  -3.81%  '
    'D43' = '5.912'
    'E43' = '  -2.79%  '
    'D44' = '1.001'
    'E44' = '  +0.17%  '
    'D45' = '7.696'
    'E45' = '  -0.52%  '
    'D46' = '0.4198'
    'E46' = '  -0.11%  '
    'D47' = '0.1248'
    'E47' = '  +0.36%  '
    'D48' = '9.040'
    'E48' = '  -1.11%  '
    'D49' = '35.07'
    'E49' = '  +0.01%  '
    'D50' = '0.05812'
    'E50' = '  -1.86%  '
    'D51' = '0.8963'
    'E51' = '  +1.30%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
